# Applies the "10-12-22 filter settings" edit:
#  1. Bumps the cached datetimeFigureOut placeholder text (master + every
#     layout) from 10/11/2022 -> 10/12/2022.
#  2. Tweaks the "Tau Chart" slide: widens/repositions the three
#     "No Reaction" callout textboxes and retitles them "No Significant
#     Bleach" (plus the tiny rounding nudge on the first callout group).
#  3. Inserts a new "Filter/Settings each Dye" slide (Title and Content
#     layout) right before the Tau Chart slide.

function EMU($emu) {
    # PowerPoint COM Left/Top/Width/Height are in points (1 pt = 12700 EMU).
    # This runtime truncates on the EMU<-pt round trip, so bias by 0.5 EMU
    # before dividing to land on the exact integer EMU value.
    return ($emu + 0.5) / 12700
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder text: 10/11/2022 -> 10/12/2022 everywhere it shows.
# ---------------------------------------------------------------------
$oldDate = "10/11/2022"
$newDate = "10/12/2022"

$master = $p.SlideMaster

for ($si = 1; $si -le $master.Shapes.Count; $si++) {
    $sh = $master.Shapes.Item($si)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    for ($si = 1; $si -le $lay.Shapes.Count; $si++) {
        $sh = $lay.Shapes.Item($si)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. "Tau Chart" slide (5th slide before the insert below) - widen the
#    three callout text boxes and update their copy.
# ---------------------------------------------------------------------
$tauChart = $p.Slides.Item(5)

# First callout group ("Group 6") needed a hairline resize too.
$grp1 = $tauChart.Shapes.Item(3)
$grp1.Left = EMU(3661963)
$grp1.Top = EMU(3977449)
$grp1.Width = EMU(2613804)
$grp1.Height = EMU(378945)

$tb1 = $grp1.GroupItems.Item(2)
$tb1.Left = EMU(2127754)
$tb1.Top = EMU(803243)
$tb1.Width = EMU(2153923)
$tb1.Height = EMU(369332)
$tb1.TextFrame.TextRange.Text = "No Significant Bleach"

# Second callout group ("Group 8") - text box only.
$grp2 = $tauChart.Shapes.Item(4)
$tb2 = $grp2.GroupItems.Item(2)
$tb2.Left = EMU(2127751)
$tb2.Top = EMU(803237)
$tb2.Width = EMU(2153923)
$tb2.Height = EMU(369332)
$tb2.TextFrame.TextRange.Text = "No Significant Bleach"

# Third callout group ("Group 11") - text box only.
$grp3 = $tauChart.Shapes.Item(5)
$tb3 = $grp3.GroupItems.Item(2)
$tb3.Left = EMU(2127751)
$tb3.Top = EMU(803237)
$tb3.Width = EMU(2153923)
$tb3.Height = EMU(369332)
$tb3.TextFrame.TextRange.Text = "No Significant Bleach"

# ---------------------------------------------------------------------
# 3. Insert the new "Filter/Settings each Dye" slide right before the
#    Tau Chart slide (position 5 -> Tau Chart becomes 6, Notes becomes 7).
# ---------------------------------------------------------------------
$titleAndContent = $master.CustomLayouts.Item(2)
$newSlide = $p.Slides.AddSlide(5, $titleAndContent)

$title = $newSlide.Shapes.Item(1)
$title.Left = EMU(191218)
$title.Top = EMU(0)
$title.Width = EMU(10515600)
$title.Height = EMU(540649)
$title.TextFrame.TextRange.Text = "Filter/Settings each Dye"

$content = $newSlide.Shapes.Item(2)
$content.Left = EMU(355120)
$content.Top = EMU(850839)
$content.Width = EMU(11583837)
$content.Height = EMU(5731115)
